$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (A1:D1) to the new "clean" column names
$ws.Cells.Item(1, 1).Value = "mx_state"
$ws.Cells.Item(1, 2).Value = "mx_municipality"
$ws.Cells.Item(1, 3).Value = "n_matriculas"
$ws.Cells.Item(1, 4).Value = "pct_matriculas"

# 2. Title-case the Spanish connector words (de/del/la/el/los/las/y) that
#    appear lower-cased inside the state (column A) and municipality
#    (column B) names, e.g. "Amatenango de la Frontera" ->
#    "Amatenango De La Frontera". Applies to data rows 2..689.
for ($r = 2; $r -le 689; $r++) {
    foreach ($col in @(1, 2)) {
        $cell = $ws.Cells.Item($r, $col)
        $v = $cell.Value2
        if ($v -ne $null) {
            $nv = $v -replace '\bde\b', 'De' `
                      -replace '\bdel\b', 'Del' `
                      -replace '\bla\b', 'La' `
                      -replace '\bel\b', 'El' `
                      -replace '\blos\b', 'Los' `
                      -replace '\blas\b', 'Las' `
                      -replace '\by\b', 'Y'
            $cell.Value = $nv
        }
    }
}

# 3. Drop the trailing metadata/footer rows (691-695): sample size,
#    source, author and date notes that no longer belong in the clean
#    dataset. This also shrinks the used range down to row 689.
$ws.Range("A691:A695").EntireRow.Delete()
